# Update the Rules sheet: cell E8 changes from "Good Morning" to "GIT UPDATE"
# (this is the real content change behind the shared-strings table churn in the diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Activate()

$ws.Range("E8").Value = "GIT UPDATE"

# Leave the active cell/selection on E8, matching the saved selection in the file
$ws.Range("E8").Select()
